$wb = $excel.ActiveWorkbook

# Sheet 1: Overview
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-03 12:18:44"
$wsOverview.Range("G3").Value = "2016-09-03 12:18:44"

# Sheet 2: zh-cn
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-09-03 12:18:39"
$wsZhCn.Range("H3").Value = "2016-09-03 12:18:39"
$wsZhCn.Range("K2").Value = "2016-09-03 12:18:55"
$wsZhCn.Range("K3").Value = "2016-09-03 12:18:55"

# Sheet 3: de-de
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-03 12:18:44"
$wsDeDe.Range("H3").Value = "2016-09-03 12:18:44"
$wsDeDe.Range("K2").Value = "2016-09-03 12:19:04"
$wsDeDe.Range("K3").Value = "2016-09-03 12:19:04"
